$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H12").Value = 6000
$ws.Range("I12").Value = 5000
$ws.Range("J12").Value = 10000
$ws.Range("K12").Value = 5000
$ws.Range("L12").Value = 10000
$ws.Range("M12").Value = -4830
$ws.Range("N12").Value = -10340

$ws.Range("H42").Value = 476.18182
$ws.Range("I42").Value = 7.6
$ws.Range("J42").Value = 866.6667
$ws.Range("K42").Value = 22.8
$ws.Range("L42").Value = 2600.0001
$ws.Range("M42").Value = 207.2
$ws.Range("N42").Value = -3060.0001

$ws.Range("H43").Value = 9249.833000000001
$ws.Range("J43").Value = 8833
$ws.Range("L43").Value = 8833
$ws.Range("N43").Value = -8971

$ws.Range("H92").Value = 578.2222
$ws.Range("I92").Value = 514.1429000000001
$ws.Range("K92").Value = 514.1429000000001
$ws.Range("M92").Value = 733.8570999999999

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 3936.5173
$ws.Range("I32").Value = 3283.4644
$ws.Range("K32").Value = 3283.4644
$ws.Range("M32").Value = -2996.4644

$ws.Range("H61").Value = 7624.75
$ws.Range("I61").Value = 7833
$ws.Range("K61").Value = 7833
$ws.Range("M61").Value = -7621

$ws.Range("H97").Value = 334.66666
$ws.Range("I97").Value = 334.66666
$ws.Range("K97").Value = 334.66666
$ws.Range("M97").Value = 161.33334

$ws.Range("H102").Value = 2230.4
$ws.Range("I102").Value = 1038.25
$ws.Range("J102").Value = 6999
$ws.Range("K102").Value = 1038.25
$ws.Range("L102").Value = 6999
$ws.Range("M102").Value = 583.75
$ws.Range("N102").Value = -10243

$ws.Range("H122").Value = 2151.1428
$ws.Range("I122").Value = 2209.6667
$ws.Range("J122").Value = 1800
$ws.Range("K122").Value = 6629.000100000001
$ws.Range("L122").Value = 5400
$ws.Range("M122").Value = -4179.000100000001
$ws.Range("N122").Value = -10300

$ws.Range("H136").Value = 7624.75
$ws.Range("I136").Value = 7833
$ws.Range("K136").Value = 23499
$ws.Range("M136").Value = -20949

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H3").Value = 5583
$ws.Range("J3").Value = 6499.6
$ws.Range("L3").Value = 6499.6
$ws.Range("N3").Value = -6725.6

$ws.Range("H7").Value = 1044.4348
$ws.Range("J7").Value = 1059.5385
$ws.Range("L7").Value = 1059.5385
$ws.Range("N7").Value = -1285.5385

$ws.Range("H10").Value = 2813.5715
$ws.Range("I10").Value = 799.5
$ws.Range("J10").Value = 5499
$ws.Range("K10").Value = 799.5
$ws.Range("L10").Value = 5499
$ws.Range("M10").Value = -660.5
$ws.Range("N10").Value = -5777

$ws.Range("H25").Value = 637.5
$ws.Range("I25").Value = 475
$ws.Range("J25").Value = 800
$ws.Range("K25").Value = 475
$ws.Range("L25").Value = 800
$ws.Range("M25").Value = -301
$ws.Range("N25").Value = -1148

$ws.Range("H99").Value = 6199.1
$ws.Range("I99").Value = 5549.923
$ws.Range("K99").Value = 5549.923
$ws.Range("M99").Value = -4051.923

$ws.Range("H122").Value = 1053
$ws.Range("I122").Value = 1053
$ws.Range("J122").Value = 0
$ws.Range("K122").Value = 3159
$ws.Range("L122").Value = 0
$ws.Range("M122").Value = -709
$ws.Range("N122").Value = ""

$ws.Range("H126").Value = 6199.1
$ws.Range("I126").Value = 5549.923
$ws.Range("K126").Value = 16649.769
$ws.Range("M126").Value = -14179.769

$ws.Range("H134").Value = 2752.5
$ws.Range("I134").Value = 2253
$ws.Range("K134").Value = 6759
$ws.Range("M134").Value = -4224

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H12").Value = 1529.7142
$ws.Range("I12").Value = 45
$ws.Range("J12").Value = 1777.1666
$ws.Range("K12").Value = 135
$ws.Range("L12").Value = 5331.4998
$ws.Range("M12").Value = 38
$ws.Range("N12").Value = -5677.4998

$ws.Range("H18").Value = 1451.6666
$ws.Range("J18").Value = 2000
$ws.Range("L18").Value = 6000
$ws.Range("N18").Value = -6338

$ws.Range("H98").Value = 600.4
$ws.Range("I98").Value = 551
$ws.Range("J98").Value = 674.5
$ws.Range("K98").Value = 1653
$ws.Range("L98").Value = 2023.5
$ws.Range("M98").Value = -155
$ws.Range("N98").Value = -5019.5

$ws.Range("H113").Value = 458
$ws.Range("J113").Value = 458
$ws.Range("L113").Value = 1374
$ws.Range("N113").Value = -5714

$ws.Range("H132").Value = 3056.9697
$ws.Range("I132").Value = 1398.6
$ws.Range("K132").Value = 12587.4
$ws.Range("M132").Value = -10057.4

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H27").Value = 0
$ws.Range("I27").Value = 0
$ws.Range("J27").Value = 0
$ws.Range("K27").Value = 0
$ws.Range("L27").Value = 0
$ws.Range("M27").Value = ""
$ws.Range("N27").Value = ""

$ws.Range("H29").Value = 10153
$ws.Range("I29").Value = 307
$ws.Range("J29").Value = 19999
$ws.Range("K29").Value = 307
$ws.Range("L29").Value = 19999
$ws.Range("M29").Value = -17
$ws.Range("N29").Value = -20579

$ws.Range("H47").Value = 30000
$ws.Range("J47").Value = 30000
$ws.Range("L47").Value = 30000
$ws.Range("N47").Value = -31136

$ws.Range("H70").Value = 9403.333000000001
$ws.Range("I70").Value = 16212
$ws.Range("J70").Value = 5999
$ws.Range("K70").Value = 16212
$ws.Range("L70").Value = 5999
$ws.Range("M70").Value = -15942
$ws.Range("N70").Value = -6539

$ws.Range("H73").Value = 9403.333000000001
$ws.Range("I73").Value = 16212
$ws.Range("J73").Value = 5999
$ws.Range("K73").Value = 16212
$ws.Range("L73").Value = 5999
$ws.Range("M73").Value = -15276
$ws.Range("N73").Value = -7871

$ws.Range("H126").Value = 7391.5
$ws.Range("I126").Value = 5725
$ws.Range("J126").Value = 8224.75
$ws.Range("K126").Value = 17175
$ws.Range("L126").Value = 24674.25
$ws.Range("M126").Value = -14705
$ws.Range("N126").Value = -29614.25

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H20").Value = 6742.5713
$ws.Range("I20").Value = 1104
$ws.Range("J20").Value = 8998
$ws.Range("K20").Value = 1104
$ws.Range("L20").Value = 8998
$ws.Range("M20").Value = -878
$ws.Range("N20").Value = -9450

$ws.Range("H22").Value = 2404.3333
$ws.Range("I22").Value = 1725.5
$ws.Range("J22").Value = 2743.75
$ws.Range("K22").Value = 1725.5
$ws.Range("L22").Value = 2743.75
$ws.Range("M22").Value = -1430.5
$ws.Range("N22").Value = -3333.75

$ws.Range("H23").Value = 0
$ws.Range("J23").Value = 0
$ws.Range("L23").Value = 0
$ws.Range("N23").Value = ""

$ws.Range("H27").Value = 2404.3333
$ws.Range("I27").Value = 1725.5
$ws.Range("J27").Value = 2743.75
$ws.Range("K27").Value = 1725.5
$ws.Range("L27").Value = 2743.75
$ws.Range("M27").Value = -1618.5
$ws.Range("N27").Value = -2957.75

$ws.Range("H40").Value = 3993.5
$ws.Range("I40").Value = 3993.3333
$ws.Range("K40").Value = 3993.3333
$ws.Range("M40").Value = -3857.3333

$ws.Range("H42").Value = 19998
$ws.Range("J42").Value = 19998
$ws.Range("L42").Value = 19998
$ws.Range("N42").Value = -21124

$ws.Range("H43").Value = 15000
$ws.Range("J43").Value = 15000
$ws.Range("L43").Value = 15000
$ws.Range("N43").Value = -15386

$ws.Range("H46").Value = 7300.4
$ws.Range("J46").Value = 5375.5
$ws.Range("L46").Value = 5375.5
$ws.Range("N46").Value = -5751.5

$ws.Range("H49").Value = 19998
$ws.Range("J49").Value = 19998
$ws.Range("L49").Value = 19998
$ws.Range("N49").Value = -20292

$ws.Range("H55").Value = 966.8333
$ws.Range("I55").Value = 0
$ws.Range("J55").Value = 966.8333
$ws.Range("K55").Value = 0
$ws.Range("L55").Value = 966.8333
$ws.Range("M55").Value = ""
$ws.Range("N55").Value = -1312.8333

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H62").Value = 1837.3334
$ws.Range("I62").Value = 1769
$ws.Range("K62").Value = 1769
$ws.Range("M62").Value = -1145

$ws.Range("H65").Value = 1837.3334
$ws.Range("I65").Value = 1769
$ws.Range("K65").Value = 8845
$ws.Range("M65").Value = -5725

$ws.Range("H70").Value = 32500
$ws.Range("I70").Value = 40000
$ws.Range("K70").Value = 40000
$ws.Range("M70").Value = -39685

$ws.Range("H73").Value = 32500
$ws.Range("I73").Value = 40000
$ws.Range("K73").Value = 40000
$ws.Range("M73").Value = -38908

$ws.Range("H76").Value = 0
$ws.Range("J76").Value = 0
$ws.Range("L76").Value = 0
$ws.Range("N76").Value = ""

$ws.Range("H79").Value = 0
$ws.Range("J79").Value = 0
$ws.Range("L79").Value = 0
$ws.Range("N79").Value = ""

$ws.Range("H122").Value = 850
$ws.Range("I122").Value = 850
$ws.Range("K122").Value = 2550
$ws.Range("M122").Value = -100

$ws.Range("H132").Value = 2199.5
$ws.Range("I132").Value = 2199.5
$ws.Range("J132").Value = 0
$ws.Range("K132").Value = 6598.5
$ws.Range("L132").Value = 0
$ws.Range("M132").Value = -4068.5
$ws.Range("N132").Value = ""

$ws.Range("H136").Value = 4207.1816
$ws.Range("J136").Value = 5199.75
$ws.Range("L136").Value = 15599.25
$ws.Range("N136").Value = -20699.25
